# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit updates to the Cerberus Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 57
$ws.Cells.Item(57, 8).Value = 89375.25
$ws.Cells.Item(57, 10).Value = 89375.25
$ws.Cells.Item(57, 12).Value = 268125.75
$ws.Cells.Item(57, 14).Value = -269123.75
# Row 132
$ws.Cells.Item(132, 8).Value = 2480.525
$ws.Cells.Item(132, 9).Value = 2295.2896
$ws.Cells.Item(132, 11).Value = 6885.8688
$ws.Cells.Item(132, 13).Value = -4355.8688
# Row 135
$ws.Cells.Item(135, 8).Value = 9198.556
$ws.Cells.Item(135, 9).Value = 2848.375
$ws.Cells.Item(135, 11).Value = 25635.375
$ws.Cells.Item(135, 13).Value = -23100.375
# Row 137
$ws.Cells.Item(137, 8).Value = 2285.2942
$ws.Cells.Item(137, 9).Value = 1004.8333
$ws.Cells.Item(137, 10).Value = 5358.4
$ws.Cells.Item(137, 11).Value = 3014.4999
$ws.Cells.Item(137, 12).Value = 16075.2
$ws.Cells.Item(137, 13).Value = -464.4998999999998
$ws.Cells.Item(137, 14).Value = -21175.2
# Row 138
$ws.Cells.Item(138, 8).Value = 3261.9014
$ws.Cells.Item(138, 10).Value = 3116.0576
$ws.Cells.Item(138, 12).Value = 9348.1728
$ws.Cells.Item(138, 14).Value = -19628.1728

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Cells.Item(132, 8).Value = 2052.65
$ws.Cells.Item(132, 9).Value = 1744.3889
$ws.Cells.Item(132, 11).Value = 5233.1667
$ws.Cells.Item(132, 13).Value = -2703.1667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 58
$ws.Cells.Item(58, 8).Value = 47797.727
$ws.Cells.Item(58, 10).Value = 47797.727
$ws.Cells.Item(58, 12).Value = 47797.727
$ws.Cells.Item(58, 14).Value = -48385.727
# Row 86
$ws.Cells.Item(86, 8).Value = 11706.214
$ws.Cells.Item(86, 9).Value = 4446.091
$ws.Cells.Item(86, 11).Value = 4446.091
$ws.Cells.Item(86, 13).Value = -3323.091
# Row 89
$ws.Cells.Item(89, 8).Value = 11706.214
$ws.Cells.Item(89, 9).Value = 4446.091
$ws.Cells.Item(89, 11).Value = 22230.455
$ws.Cells.Item(89, 13).Value = -16614.455
# Row 96
$ws.Cells.Item(96, 8).Value = 12248.75
$ws.Cells.Item(96, 9).Value = 12248.75
$ws.Cells.Item(96, 11).Value = 12248.75
$ws.Cells.Item(96, 13).Value = -9502.75
# Row 105
$ws.Cells.Item(105, 8).Value = 2850.4443
$ws.Cells.Item(105, 9).Value = 2850.4443
$ws.Cells.Item(105, 11).Value = 2850.4443
$ws.Cells.Item(105, 13).Value = -1103.4443
# Row 134
$ws.Cells.Item(134, 8).Value = 7572.115
$ws.Cells.Item(134, 9).Value = 7325.033
$ws.Cells.Item(134, 11).Value = 21975.099
$ws.Cells.Item(134, 13).Value = -19440.099

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 94.23077000000001
$ws.Cells.Item(7, 9).Value = 101.625
$ws.Cells.Item(7, 10).Value = 82.40000000000001
$ws.Cells.Item(7, 11).Value = 101.625
$ws.Cells.Item(7, 12).Value = 82.40000000000001
$ws.Cells.Item(7, 13).Value = 11.375
$ws.Cells.Item(7, 14).Value = -308.4
# Row 19
$ws.Cells.Item(19, 8).Value = 1178084.4
$ws.Cells.Item(19, 9).Value = 1668027.1
$ws.Cells.Item(19, 10).Value = 2222
$ws.Cells.Item(19, 11).Value = 1668027.1
$ws.Cells.Item(19, 12).Value = 2222
$ws.Cells.Item(19, 13).Value = -1667857.1
$ws.Cells.Item(19, 14).Value = -2562
# Row 24
$ws.Cells.Item(24, 8).Value = 1178084.4
$ws.Cells.Item(24, 9).Value = 1668027.1
$ws.Cells.Item(24, 10).Value = 2222
$ws.Cells.Item(24, 11).Value = 1668027.1
$ws.Cells.Item(24, 12).Value = 2222
$ws.Cells.Item(24, 13).Value = -1667857.1
$ws.Cells.Item(24, 14).Value = -2562
# Row 31
$ws.Cells.Item(31, 8).Value = 1704.875
$ws.Cells.Item(31, 9).Value = 874.2941
$ws.Cells.Item(31, 10).Value = 2646.2
$ws.Cells.Item(31, 11).Value = 874.2941
$ws.Cells.Item(31, 12).Value = 2646.2
$ws.Cells.Item(31, 13).Value = -579.2941
$ws.Cells.Item(31, 14).Value = -3236.2
# Row 34
$ws.Cells.Item(34, 8).Value = 1704.875
$ws.Cells.Item(34, 9).Value = 874.2941
$ws.Cells.Item(34, 10).Value = 2646.2
$ws.Cells.Item(34, 11).Value = 874.2941
$ws.Cells.Item(34, 12).Value = 2646.2
$ws.Cells.Item(34, 13).Value = -672.2941
$ws.Cells.Item(34, 14).Value = -3050.2
# Row 107
$ws.Cells.Item(107, 8).Value = 2918.5
$ws.Cells.Item(107, 9).Value = 3192.7222
$ws.Cells.Item(107, 11).Value = 3192.7222
$ws.Cells.Item(107, 13).Value = -1272.7222

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Cells.Item(39, 8).Value = 7999.2
$ws.Cells.Item(39, 10).Value = 7776.8887
$ws.Cells.Item(39, 12).Value = 23330.6661
$ws.Cells.Item(39, 14).Value = -23918.6661
# Row 103
$ws.Cells.Item(103, 8).Value = 5446.778
$ws.Cells.Item(103, 9).Value = 4512
$ws.Cells.Item(103, 10).Value = 5713.857
$ws.Cells.Item(103, 11).Value = 13536
$ws.Cells.Item(103, 12).Value = 17141.571
$ws.Cells.Item(103, 13).Value = -12657
$ws.Cells.Item(103, 14).Value = -18899.571
# Row 105
$ws.Cells.Item(105, 8).Value = 9333
$ws.Cells.Item(105, 10).Value = 9333
$ws.Cells.Item(105, 12).Value = 27999
$ws.Cells.Item(105, 14).Value = -33241
# Row 107
$ws.Cells.Item(107, 8).Value = 1199.5122
$ws.Cells.Item(107, 9).Value = 240.75
$ws.Cells.Item(107, 10).Value = 1431.9395
$ws.Cells.Item(107, 11).Value = 722.25
$ws.Cells.Item(107, 12).Value = 4295.818499999999
$ws.Cells.Item(107, 13).Value = 1197.75
$ws.Cells.Item(107, 14).Value = -8135.818499999999
# Row 111
$ws.Cells.Item(111, 8).Value = 3377.2
$ws.Cells.Item(111, 9).Value = 1471.75
$ws.Cells.Item(111, 11).Value = 4415.25
$ws.Cells.Item(111, 13).Value = -1348.25
# Row 117
$ws.Cells.Item(117, 8).Value = 3829.1333
$ws.Cells.Item(117, 9).Value = 2633.3333
$ws.Cells.Item(117, 10).Value = 4128.0835
$ws.Cells.Item(117, 11).Value = 7899.999899999999
$ws.Cells.Item(117, 12).Value = 12384.2505
$ws.Cells.Item(117, 13).Value = -4457.999899999999
$ws.Cells.Item(117, 14).Value = -19268.2505
# Row 118
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 13).ClearContents()
# Row 119
$ws.Cells.Item(119, 8).Value = 2998
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
# Row 120
$ws.Cells.Item(120, 8).Value = 39999.855
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 39999.855
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).ClearContents()
$ws.Cells.Item(120, 13).Value = 119999.565
$ws.Cells.Item(120, 14).Value = -129675.565
# Row 121
$ws.Cells.Item(121, 8).Value = 76360.44500000001
$ws.Cells.Item(121, 9).Value = 25522.25
$ws.Cells.Item(121, 10).Value = 90885.64
$ws.Cells.Item(121, 11).Value = 76566.75
$ws.Cells.Item(121, 12).Value = 272656.92
$ws.Cells.Item(121, 13).Value = -75256.75
$ws.Cells.Item(121, 14).Value = -275276.92
# Row 137
$ws.Cells.Item(137, 8).Value = 23492.25
$ws.Cells.Item(137, 9).Value = 3969
$ws.Cells.Item(137, 11).Value = 11907
$ws.Cells.Item(137, 13).Value = -6807
# Row 141
$ws.Cells.Item(141, 8).Value = 42879.707
$ws.Cells.Item(141, 9).Value = 9651.666999999999
$ws.Cells.Item(141, 11).Value = 28955.001
$ws.Cells.Item(141, 13).Value = -23775.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 28998.666
$ws.Cells.Item(15, 10).Value = 28998.666
$ws.Cells.Item(15, 12).Value = 28998.666
$ws.Cells.Item(15, 14).Value = -29574.666
# Row 81
$ws.Cells.Item(81, 8).Value = 28998.666
$ws.Cells.Item(81, 10).Value = 28998.666
$ws.Cells.Item(81, 12).Value = 28998.666
$ws.Cells.Item(81, 14).Value = -30994.666
# Row 84
$ws.Cells.Item(84, 8).Value = 28998.666
$ws.Cells.Item(84, 10).Value = 28998.666
$ws.Cells.Item(84, 12).Value = 86995.99800000001
$ws.Cells.Item(84, 14).Value = -96979.99800000001
# Row 113
$ws.Cells.Item(113, 8).Value = 2099.0833
$ws.Cells.Item(113, 9).Value = 1649.25
$ws.Cells.Item(113, 11).Value = 1649.25
$ws.Cells.Item(113, 13).Value = 520.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 2090.4
$ws.Cells.Item(46, 9).Value = 911
$ws.Cells.Item(46, 11).Value = 911
$ws.Cells.Item(46, 13).Value = -723
# Row 82
$ws.Cells.Item(82, 8).Value = 2004.0869
$ws.Cells.Item(82, 9).Value = 1840.875
$ws.Cells.Item(82, 10).Value = 2377.1428
$ws.Cells.Item(82, 11).Value = 1840.875
$ws.Cells.Item(82, 12).Value = 2377.1428
$ws.Cells.Item(82, 13).Value = -1479.875
$ws.Cells.Item(82, 14).Value = -3099.1428
# Row 85
$ws.Cells.Item(85, 8).Value = 2004.0869
$ws.Cells.Item(85, 9).Value = 1840.875
$ws.Cells.Item(85, 10).Value = 2377.1428
$ws.Cells.Item(85, 11).Value = 1840.875
$ws.Cells.Item(85, 12).Value = 2377.1428
$ws.Cells.Item(85, 13).Value = -592.875
$ws.Cells.Item(85, 14).Value = -4873.1428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 952.6667
$ws.Cells.Item(107, 9).Value = 1039.2
$ws.Cells.Item(107, 10).Value = 844.5
$ws.Cells.Item(107, 11).Value = 3117.6
$ws.Cells.Item(107, 12).Value = 2533.5
$ws.Cells.Item(107, 13).Value = -1197.6
$ws.Cells.Item(107, 14).Value = -6373.5
